# update student from fill up class 6
#
# Starting state: a single sheet "gunjon" with one student's admission info
# in column B (labels in column A). This edit duplicates that sheet four
# times to hold four more students (anjoli, bisal, Ria moni, Ripon), fills
# in three of them, leaves "bisal" mostly blank, colors the new tabs, and
# reorders the tabs to: anjoli, bisal, Ria moni, Ripon, gunjon.

$wb = $excel.ActiveWorkbook
$gunjon = $wb.Worksheets.Item("gunjon")

# --- 1. Duplicate "gunjon" to create the four new sheets -------------------
# Each Copy($null, $gunjon) places the clone immediately after "gunjon",
# so creating them in this order (Ria moni, bisal, Ripon, anjoli) yields
# sheetId 5, 6, 7, 8 respectively, and leaves them sitting (in the tab
# strip) as: gunjon, anjoli, Ripon, bisal, Ria moni  -- we fix the order later.
$gunjon.Copy($null, $gunjon)
$wb.Worksheets.Item("gunjon (2)").Name = "Ria moni"

$gunjon.Copy($null, $gunjon)
$wb.Worksheets.Item("gunjon (2)").Name = "bisal"

$gunjon.Copy($null, $gunjon)
$wb.Worksheets.Item("gunjon (2)").Name = "Ripon"

$gunjon.Copy($null, $gunjon)
$wb.Worksheets.Item("gunjon (2)").Name = "anjoli"

# --- 2. Fill in student data (column B) ------------------------------------
# Row layout (same on every sheet):
#  1 নাম                     (name)
#  2 জন্ম তারিখ               (date of birth)
#  3 জন্ম নিবন্ধনের নম্বর       (birth reg. number)
#  4 মোবাইল নম্বর              (mobile number)
#  5 পিতার নাম                (father's name)
#  6 পিতার এনআইডি নম্বর        (father's NID)
#  7 মাতার নাম                (mother's name)
#  8 মাতার এনআইডি নম্বর        (mother's NID)
#  9 পূর্ববর্তী স্কুলের নাম      (previous school)
# 10 বিস্তারিত ঠিকানা          (address)

# Fill order matters for shared-string layout: Ripon, then Ria moni, then
# anjoli (bisal is left incomplete, matching the source edit).
$ripon = $wb.Worksheets.Item("Ripon")
$ripon.Range("B1").Value = "Ripon Chandro Roy"
$ripon.Range("B2").Value = "31-05-2015"
$ripon.Range("B3").Value = "20157316431043106"
$ripon.Range("B4").Value = "01873396263"
$ripon.Range("B5").Value = "Sumanth Chandro roy"
$ripon.Range("B6").Value = "9562222860"
$ripon.Range("B7").Value = "Sumitra"
$ripon.Range("B8").Value = "5524450003"
$ripon.Range("B9").Value = "Prone Community Primary School"
$ripon.Range("B10").Value = "Kanial khata baniya para"

$riamoni = $wb.Worksheets.Item("Ria moni")
$riamoni.Range("B1").Value = "Mst Riya Moni"
$riamoni.Range("B2").Value = "28-03-2013"
$riamoni.Range("B3").Value = "20137316431038599"
$riamoni.Range("B4").Value = "01737082418"
$riamoni.Range("B5").Value = "Md. Azinur Rahman"
$riamoni.Range("B6").Value = "1907648040"
$riamoni.Range("B7").Value = "sarmin"
$riamoni.Range("B8").Value = "6460123679"
$riamoni.Range("B9").Value = "Prone Community Primary School"
$riamoni.Range("B10").Value = "Kanial khata "

$anjoli = $wb.Worksheets.Item("anjoli")
$anjoli.Range("B1").Value = "Sreemoti Anjoli Rani ray Khushi"
$anjoli.Range("B2").Value = "15-10-2015"
$anjoli.Range("B3").Value = "20157316431036135"
$anjoli.Range("B4").Value = "01785653041"
$anjoli.Range("B5").Value = "bimol chandro"
$anjoli.Range("B6").Value = "19947316431000244"
$anjoli.Range("B7").Value = "Sreemoti dolly Rani ray "
$anjoli.Range("B8").Value = "1491503601"
$anjoli.Range("B9").Value = "Prone Community Primary School"
$anjoli.Range("B10").Value = "Kanial khata"

# "bisal" keeps its blank rows 1-8 (cleared below) but rows 9-10 are filled.
$bisal = $wb.Worksheets.Item("bisal")
$bisal.Range("B1").ClearContents()
$bisal.Range("B2").ClearContents()
$bisal.Range("B3").ClearContents()
$bisal.Range("B4").ClearContents()
$bisal.Range("B5").ClearContents()
$bisal.Range("B6").ClearContents()
$bisal.Range("B7").ClearContents()
$bisal.Range("B8").ClearContents()
$bisal.Range("B9").Value = "Prone Community Primary School"
$bisal.Range("B10").Value = "Kanial khata mastar para"

# --- 3. Tab colors -----------------------------------------------------
# anjoli / bisal / Ria moni use the light/white "theme 0" tab color;
# Ripon keeps the same green as gunjon.
$anjoli.Tab.Color = 16777215
$bisal.Tab.Color = 16777215
$riamoni.Tab.Color = 16777215
$ripon.Tab.Color = 5287936

# --- 4. Reorder tabs: anjoli, bisal, Ria moni, Ripon, gunjon ----------------
$wb.Worksheets.Item("anjoli").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("bisal").Move($null, $wb.Worksheets.Item("anjoli"))
$wb.Worksheets.Item("Ria moni").Move($null, $wb.Worksheets.Item("bisal"))
$wb.Worksheets.Item("Ripon").Move($null, $wb.Worksheets.Item("Ria moni"))

# --- 5. Selection / active sheet -------------------------------------------
# gunjon used to be the active sheet (tabSelected, cell D17 selected);
# anjoli is now the active sheet, selection on each new sheet sits at B10
# (B14 for the still-being-filled-in "bisal").
$wb.Worksheets.Item("anjoli").Activate()
$wb.Worksheets.Item("anjoli").Range("B10").Select()
$wb.Worksheets.Item("bisal").Range("B14").Select()
$wb.Worksheets.Item("Ria moni").Range("B10").Select()
$wb.Worksheets.Item("Ripon").Range("B10").Select()

$wb.Worksheets.Item("anjoli").Activate()
